$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.676.87'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.596.01'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '1.819.86'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.02'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.620.58'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.01'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '26.642.26'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '0.0₃0751'
$ws.Range('E18').Value = '  +3.08%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '208.87'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.97'
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.31'
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.31'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.95'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').Value = '1.283.05'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.615'
$ws.Range('E35').Value = '  -7.04%  '
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.05'
$ws.Range('E39').Value = '  +16.23%  '
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.90'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').Value = '1.732.77'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.62'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0508'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.30'
$ws.Range('E51').Value = '  -2.63%  '
